$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the row that only carried the "Docentes responsaveis" value
# in columns B/C with an empty label in column A). This shifts rows 14-22 up
# by one, matching the row count/height layout of the edited workbook.
$ws.Rows.Item(13).Delete()

# After the shift, fix up the handful of cells whose text content changed
# (per the authoritative diff) so they no longer just mirror the row that
# used to sit below them.
$ws.Range("B10").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C10").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("B15").Value = '01/01/2022'
$ws.Range("C15").Value = '01/01/2022'
$ws.Range("B18").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C18").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("B19").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Range("C19").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Range("B20").Value = 'MF = média finalMF = (P1 + P2)/2'
$ws.Range("C20").Value = 'MF = média finalMF = (P1 + P2)/2'
$ws.Range("B21").Value = 'Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Range("C21").Value = 'Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'

Write-Host "Applied LOT2045 content update"
